# The document contains four ${...} placeholder tokens whose internal
# variable names need to be renamed to the new "customers table" field
# names:
#   ugyfelnev    -> customername
#   cim          -> address
#   vezeto       -> leader
#   telefonszam  -> phone
#
# Use Find/Replace (whole-word, case-sensitive) scoped to the whole
# document body so only the placeholder identifiers themselves are
# touched, leaving the surrounding "${" / "}" text and layout untouched.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "ugyfelnev";   New = "customername" },
    @{ Old = "cim";         New = "address" },
    @{ Old = "vezeto";      New = "leader" },
    @{ Old = "telefonszam"; New = "phone" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.New, 2) | Out-Null
}
